$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + shifted/updated edge list (topological similarity reordering)
$data = @(
    @("source", "Target"),
    @("A", "C"),
    @("A", "D"),
    @("A", "F"),
    @("A", "M"),
    @("A", "B"),
    @("A", "G"),
    @("B", "M"),
    @("B", "G"),
    @("B", "F"),
    @("B", "Q"),
    @("Q", "I"),
    @("Q", "F"),
    @("Q", "W"),
    @("C", "I"),
    @("C", "G")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("D6").Select()
